$d = $word.ActiveDocument

function Escape-Xml($text) {
    $t = $text -replace "&", "&amp;"
    $t = $t -replace "<", "&lt;"
    $t = $t -replace ">", "&gt;"
    return $t
}

function Find-ParagraphByText($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $needle) {
            return $p
        }
    }
    return $null
}

# Replaces the *entire* content of a paragraph with a single run holding
# $newText (optionally wrapped in $rPrXml formatting, e.g. "<w:b/>").
# Because the replacement range begins exactly at the paragraph's Start,
# any pre-existing zero-length "<w:r/>" marker run just before it is left
# untouched, matching how this document structures its runs.
function Set-ParagraphRunText($paragraph, $newText, $rPrXml) {
    $rng = $d.Range($paragraph.Range.Start, $paragraph.Range.End)
    $escaped = Escape-Xml $newText
    if ($rPrXml) {
        $runXml = '<w:r><w:rPr>' + $rPrXml + '</w:rPr><w:t>' + $escaped + '</w:t></w:r>'
    } else {
        $runXml = '<w:r><w:t>' + $escaped + '</w:t></w:r>'
    }
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

function Insert-BulletAfter($paragraph, $text) {
    $paragraph.Range.InsertParagraphAfter()
    $newPara = $paragraph.Next()
    $escaped = Escape-Xml $text
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>' + $escaped + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $newPara.Range.InsertXML($xml)
    return $newPara
}

# --- Simple text replacements (titles / headings / bold / list items) ---
$p = Find-ParagraphByText "Play Grand Spinn Slot Game for Free"
Set-ParagraphRunText $p "Play Grand Spinn for Free" $null

$p = Find-ParagraphByText "Simple and linear gameplay mechanics"
Set-ParagraphRunText $p "Traditional gameplay experience" $null

$p = Find-ParagraphByText "Only features a 3x3 grid and single payline"
Set-ParagraphRunText $p "Limited number of symbols per play" $null

$p = Find-ParagraphByText "Not suitable for those who prefer more complex gameplay"
Set-ParagraphRunText $p "Only one payline" $null

$p = Find-ParagraphByText "Play Grand Spinn Slot Game for Free"
Set-ParagraphRunText $p "Play Grand Spinn for Free" "<w:b/>"

$p = Find-ParagraphByText "Read our review of Grand Spinn, a traditional slot game with Wild symbols and three Jackpots. Play for free and enjoy its nostalgic graphics and design."
Set-ParagraphRunText $p "Read our review of Grand Spinn and play this engaging game for free." "<w:i/>"

# --- Restructure the "What we like" bullet list ---
# Insert the two new bullet points after "Traditional gameplay experience"
$anchor = Find-ParagraphByText "Traditional gameplay experience"
$anchor = Insert-BulletAfter $anchor "Beautifully illustrated symbols"
Insert-BulletAfter $anchor "Inclusion of Wild symbols and Jackpots" | Out-Null

# Remove the two obsolete bullet points
$p = Find-ParagraphByText "Wild symbols and three different Jackpots"
if ($p -ne $null) { $p.Range.Delete() | Out-Null }

$p = Find-ParagraphByText "Appealing to a wide audience"
if ($p -ne $null) { $p.Range.Delete() | Out-Null }
